# "fixed problem with headless mode"
# Adds a new patient row (MUSC 312-B / left) beneath the existing data
# and moves the active selection down to the next empty row, the way a
# user would after typing a new entry into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "MUSC 312-B"
$ws.Range("B3").Value = "left"

$ws.Range("A4").Select() | Out-Null
